# [PV-350][WIP] Replace hard coding of visual height with calculated value
#
# Update the header row of the "PV-Test-03" sheet to use the new column
# header labels:
#   "Row ID"     -> "Id"
#   "Task"       -> "Task Name"
#   "Start Date" -> "Start"
#   "End Date"   -> "Finish"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03")

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"
